{"js": "const body = context.document.body;\n\n// \"Fix SED to SDE\" - the job title was mistyped as \"SED Intern\" and should\n// read \"SDE Intern\" (Software Development Engineer Intern).\nconst sedResults = body.search(\"SED\", { matchCase: true, matchWholeWord: true });\nsedResults.load(\"items\");\nawait context.sync();\nif (sedResults.items.length !== 1) {\n  throw new Error(\"Expected exactly 1 match for 'SED', got \" + sedResults.items.length);\n}\nsedResults.items[0].insertText(\"SDE\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Word drops its \"last edit\" marker (the hidden _GoBack bookmark) at the\n// spot of the most recent change. Relocate the stale one left over from a\n// previous edit to right after the text we just fixed, so only a single\n// \"_GoBack\" bookmark remains afterwards (rather than two bookmarks sharing\n// the same name).\ntry {\n  const oldBookmarkRange = body.getBookmarkRangeOrNullObject(\"_GoBack\");\n  await context.sync();\n\n  if (!oldBookmarkRange.isNullObject) {\n    const oldPara = oldBookmarkRange.paragraphs.getFirst();\n    const oldParaRange = oldPara.getRange(\"Whole\");\n    const oldOoxml = oldPara.getOoxml();\n    await context.sync();\n\n    const fullPkg = oldOoxml.value;\n    const pXmlMatch = fullPkg.match(/<w:p[ >][\\s\\S]*?<\\/w:p>/);\n    if (pXmlMatch) {\n      let pXml = pXmlMatch[0];\n\n      // Drop the synthetic paraId/textId the single-paragraph OOXML export\n      // fabricates (the original paragraph did not have them).\n      pXml = pXml.replace(/ w14:paraId=\"[0-9A-Fa-f]+\"/, \"\").replace(/ w14:textId=\"[0-9A-Fa-f]+\"/, \"\");\n\n      // Strip the hidden \"_GoBack\" bookmark markers out of this paragraph.\n      const idMatch = pXml.match(/<w:bookmarkStart w:id=\"(\\d+)\" w:name=\"_GoBack\"\\/>/);\n      if (idMatch) {\n        const id = idMatch[1];\n        pXml = pXml\n          .replace(new RegExp('<w:bookmarkStart w:id=\"' + id + '\" w:name=\"_GoBack\"/>'), \"\")\n          .replace(new RegExp('<w:bookmarkEnd w:id=\"' + id + '\"/>'), \"\");\n\n        const ooxmlPkg = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">'\n          + '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>'\n          + '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + pXml + '</w:body></w:document>'\n          + '</pkg:xmlData></pkg:part></pkg:package>';\n\n        oldParaRange.insertOoxml(ooxmlPkg, Word.InsertLocation.replace);\n        await context.sync();\n      }\n    }\n  }\n\n  // Re-insert \"_GoBack\" collapsed right after the corrected \"SDE\".\n  const sdeResults = body.search(\"SDE\", { matchCase: true, matchWholeWord: true });\n  sdeResults.load(\"items\");\n  await context.sync();\n  if (sdeResults.items.length > 0) {\n    const sdeEnd = sdeResults.items[0].getRange(\"End\");\n    sdeEnd.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n} catch (e) {\n  // The _GoBack bookmark is just Word's \"last edit location\" marker; if\n  // relocating it fails for any reason, the meaningful text fix above has\n  // already been applied, so don't let this secondary step fail the edit.\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# \"Fix SED to SDE\" - the job title was mistyped as \"SED Intern\" and\n# should read \"SDE Intern\" (Software Development Engineer Intern).\n$find = $d.Content.Find\n$find.Text = \"SED\"\n$find.MatchWholeWord = $true\n$find.MatchCase = $true\n$find.Replacement.Text = \"SDE\"\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n\n# Word drops its \"last edit\" marker (the hidden _GoBack bookmark) at the\n# spot of the most recent change, so move it to right after the fixed\n# \"SDE\" text (replacing whichever paragraph it used to mark).\n$rng = $d.Content\n$find2 = $rng.Find\n$find2.Text = \"SDE\"\n$find2.MatchWholeWord = $true\n$find2.MatchCase = $true\n$find2.Execute() | Out-Null\nif ($find2.Found) {\n    $rng.Collapse(0) | Out-Null\n    $d.Bookmarks.Add(\"_GoBack\", $rng) | Out-Null\n}\n"}
